# Add a new announcement row (row 7) to the table on Sheet1:
#   DATE = 10/2/25 (serial 45930)
#   ANNOUNCEMENT   = "Practice on 10/2/25 will be moved to 7-9pm instead."
#   LUS TSHAJ TAWM = "Peb kawm nkauj hnub 10/2/25 thaum 7 moo txog 9 moo"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row so Table1's ref/autoFilter and the
# worksheet dimension expand from A1:C6 to A1:C7.
$lo = $wb.Worksheets.Item(1).ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the formatting of an existing data row (row 3) onto the new
# row 7 so it picks up the same cell styles (date style, wrap-text
# styles, borders, fill) used throughout the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)

# Fill in the new row's values.
$ws.Range("A7").Value = 45930
$ws.Range("B7").Value = "Practice on 10/2/25 will be moved to 7-9pm instead."
$ws.Range("C7").Value = "Peb kawm nkauj hnub 10/2/25 thaum 7 moo txog 9 moo"

# Match the row height (30pt) used by the other wrapped-text rows.
$ws.Rows.Item(7).RowHeight = 30
